$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update section title wording: "MWF" -> "M-W-F"
$ws.Range("A1").Value = "CS 320 Section 102 (M-W-F 11:00 - 11:50)"

# Roster corrections
# Row 9: Kettula - first name Jordan -> John
$ws.Range("B9").Value = "John"

# Row 11: Mcbride -> McBride (capitalization fix)
$ws.Range("A11").Value = "McBride"

# Row 12: Mccloskey - first name James -> Trey
$ws.Range("B12").Value = "Trey"

# Update selection to the header row, matching the saved view state
$ws.Range("A1:E1").Select() | Out-Null
